# Generate Report for Handoff
# The e2e/6489d606-5d3b-409c-bc93-e42c01cf99eb.md file has moved from
# "handed back" to "ready for handoff" again because the handback version
# is stale. Update the status/report rows for that file across all three
# sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fd90ea9b466ce18464b4b13f46e34273ca51b360/e2e/6489d606-5d3b-409c-bc93-e42c01cf99eb.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7289950c27bd671277e6b06a673f541ec2cc68f7/e2e/6489d606-5d3b-409c-bc93-e42c01cf99eb.md."

# --- Overview sheet: row 3 is the 6489d606-... file ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-01 06:55:25"

# --- zh-cn sheet: row 3 is the 6489d606-... file ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-09-01 06:55:20"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Range("P1").ColumnWidth = 39.1666666667

# --- de-de sheet: row 3 is the 6489d606-... file ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-09-01 06:55:25"
$dede.Range("P3").Value = $errorDetail
$dede.Range("P1").ColumnWidth = 39.1666666667
